$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.06500340741874241
$ws.Range("J2").Value = 0.06500340741874241
$ws.Range("M2").Value = 0.5396133333333334
$ws.Range("N2").Value = 1.61884
$ws.Range("O2").Value = 0.02587018426425635
$ws.Range("P2").Value = 0.02587018426425635
$ws.Range("Q2").Value = 0.02628456546666667
$ws.Range("R2").Value = 0.2365610892
$ws.Range("S2").Value = 0.001681650127727395
$ws.Range("T2").Value = 0.001681650127727395
$ws.Range("I3").Value = 0.06500340741874241
$ws.Range("J3").Value = 0.06500340741874241
$ws.Range("M3").Value = 6.260434
$ws.Range("N3").Value = 18.781302
$ws.Range("O3").Value = 0.3001382122153186
$ws.Range("P3").Value = 0.3001382122153186
$ws.Range("Q3").Value = 0.30494574014
$ws.Range("R3").Value = 2.74451166126
$ws.Range("S3").Value = 0.01951000649056533
$ws.Range("T3").Value = 0.01951000649056533
$ws.Range("I4").Value = 0.06500340741874241
$ws.Range("J4").Value = 0.06500340741874241
$ws.Range("M4").Value = 1.182122
$ws.Range("N4").Value = 3.546366
$ws.Range("O4").Value = 0.0566733845769154
$ws.Range("P4").Value = 0.0566733845769154
$ws.Range("Q4").Value = 0.05758116262
$ws.Range("R4").Value = 0.51823046358
$ws.Range("S4").Value = 0.003683963107452305
$ws.Range("T4").Value = 0.003683963107452305
$ws.Range("I5").Value = 0.06500340741874241
$ws.Range("J5").Value = 0.06500340741874241
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 12.87633433333333
$ws.Range("N5").Value = 38.629003
$ws.Range("O5").Value = 0.6173182189435098
$ws.Range("P5").Value = 0.6173182189435098
$ws.Range("Q5").Value = 0.6272062453766667
$ws.Range("R5").Value = 5.64485620839
$ws.Range("S5").Value = 0.0401277876929974
$ws.Range("T5").Value = 0.0401277876929974
$ws.Range("I6").Value = 0.07892489266186128
$ws.Range("J6").Value = 0.07892489266186128
$ws.Range("M6").Value = 0.5396133333333334
$ws.Range("N6").Value = 1.61884
$ws.Range("O6").Value = 0.02587018426425635
$ws.Range("P6").Value = 0.02587018426425635
$ws.Range("Q6").Value = 0.03191381176
$ws.Range("R6").Value = 0.28722430584
$ws.Range("S6").Value = 0.002041801516199006
$ws.Range("T6").Value = 0.002041801516199006
$ws.Range("I7").Value = 0.07892489266186128
$ws.Range("J7").Value = 0.07892489266186128
$ws.Range("M7").Value = 6.260434
$ws.Range("N7").Value = 18.781302
$ws.Range("O7").Value = 0.3001382122153186
$ws.Range("P7").Value = 0.3001382122153186
$ws.Range("Q7").Value = 0.370254587628
$ws.Range("R7").Value = 3.332291288652
$ws.Range("S7").Value = 0.02368837618281696
$ws.Range("T7").Value = 0.02368837618281696
$ws.Range("I8").Value = 0.07892489266186128
$ws.Range("J8").Value = 0.07892489266186128
$ws.Range("M8").Value = 1.182122
$ws.Range("N8").Value = 3.546366
$ws.Range("O8").Value = 0.0566733845769154
$ws.Range("P8").Value = 0.0566733845769154
$ws.Range("Q8").Value = 0.06991305932399999
$ws.Range("R8").Value = 0.629217533916
$ws.Range("S8").Value = 0.004472940794517434
$ws.Range("T8").Value = 0.004472940794517434
$ws.Range("I9").Value = 0.07892489266186128
$ws.Range("J9").Value = 0.07892489266186128
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 12.87633433333333
$ws.Range("N9").Value = 38.629003
$ws.Range("O9").Value = 0.6173182189435098
$ws.Range("P9").Value = 0.6173182189435098
$ws.Range("Q9").Value = 0.761532165142
$ws.Range("R9").Value = 6.853789486277999
$ws.Range("S9").Value = 0.0487217741683279
$ws.Range("T9").Value = 0.0487217741683279
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.3835156666666666
$ws.Range("H10").Value = 1.150547
$ws.Range("I10").Value = 0.511800967600163
$ws.Range("J10").Value = 0.511800967600163
$ws.Range("M10").Value = 0.5396133333333334
$ws.Range("N10").Value = 1.61884
$ws.Range("O10").Value = 0.02587018426425635
$ws.Range("P10").Value = 0.02587018426425635
$ws.Range("Q10").Value = 0.2069501672755556
$ws.Range("R10").Value = 1.86255150548
$ws.Range("S10").Value = 0.01324038533844091
$ws.Range("T10").Value = 0.01324038533844091
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.3835156666666666
$ws.Range("H11").Value = 1.150547
$ws.Range("I11").Value = 0.511800967600163
$ws.Range("J11").Value = 0.511800967600163
$ws.Range("M11").Value = 6.260434
$ws.Range("N11").Value = 18.781302
$ws.Range("O11").Value = 0.3001382122153186
$ws.Range("P11").Value = 0.3001382122153186
$ws.Range("Q11").Value = 2.400974519132666
$ws.Range("R11").Value = 21.608770672194
$ws.Range("S11").Value = 0.1536110274255831
$ws.Range("T11").Value = 0.1536110274255831
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.3835156666666666
$ws.Range("H12").Value = 1.150547
$ws.Range("I12").Value = 0.511800967600163
$ws.Range("J12").Value = 0.511800967600163
$ws.Range("M12").Value = 1.182122
$ws.Range("N12").Value = 3.546366
$ws.Range("O12").Value = 0.0566733845769154
$ws.Range("P12").Value = 0.0566733845769154
$ws.Range("Q12").Value = 0.4533623069113333
$ws.Range("R12").Value = 4.080260762202
$ws.Range("S12").Value = 0.02900549306364146
$ws.Range("T12").Value = 0.02900549306364146
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.3835156666666666
$ws.Range("H13").Value = 1.150547
$ws.Range("I13").Value = 0.511800967600163
$ws.Range("J13").Value = 0.511800967600163
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 12.87633433333333
$ws.Range("N13").Value = 38.629003
$ws.Range("O13").Value = 0.6173182189435098
$ws.Range("P13").Value = 0.6173182189435098
$ws.Range("Q13").Value = 4.938275946071221
$ws.Range("R13").Value = 44.44448351464099
$ws.Range("S13").Value = 0.3159440617724976
$ws.Range("T13").Value = 0.3159440617724976
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.1956336666666667
$ws.Range("H14").Value = 0.586901
$ws.Range("I14").Value = 0.2610727764146126
$ws.Range("J14").Value = 0.2610727764146126
$ws.Range("M14").Value = 0.5396133333333334
$ws.Range("N14").Value = 1.61884
$ws.Range("O14").Value = 0.02587018426425635
$ws.Range("P14").Value = 0.02587018426425635
$ws.Range("Q14").Value = 0.1055665349822222
$ws.Range("R14").Value = 0.95009881484
$ws.Range("S14").Value = 0.006754000832227029
$ws.Range("T14").Value = 0.006754000832227027
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.1956336666666667
$ws.Range("H15").Value = 0.586901
$ws.Range("I15").Value = 0.2610727764146126
$ws.Range("J15").Value = 0.2610727764146126
$ws.Range("M15").Value = 6.260434
$ws.Range("N15").Value = 18.781302
$ws.Range("O15").Value = 0.3001382122153186
$ws.Range("P15").Value = 0.3001382122153186
$ws.Range("Q15").Value = 1.224751658344667
$ws.Range("R15").Value = 11.022764925102
$ws.Range("S15").Value = 0.07835791637117143
$ws.Range("T15").Value = 0.07835791637117141
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.1956336666666667
$ws.Range("H16").Value = 0.586901
$ws.Range("I16").Value = 0.2610727764146126
$ws.Range("J16").Value = 0.2610727764146126
$ws.Range("M16").Value = 1.182122
$ws.Range("N16").Value = 3.546366
$ws.Range("O16").Value = 0.0566733845769154
$ws.Range("P16").Value = 0.0566733845769154
$ws.Range("Q16").Value = 0.2312628613073333
$ws.Range("R16").Value = 2.081365751766
$ws.Range("S16").Value = 0.01479587786030839
$ws.Range("T16").Value = 0.01479587786030839
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.1956336666666667
$ws.Range("H17").Value = 0.586901
$ws.Range("I17").Value = 0.2610727764146126
$ws.Range("J17").Value = 0.2610727764146126
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 12.87633433333333
$ws.Range("N17").Value = 38.629003
$ws.Range("O17").Value = 0.6173182189435098
$ws.Range("P17").Value = 0.6173182189435098
$ws.Range("Q17").Value = 2.519044498855889
$ws.Range("R17").Value = 22.671400489703
$ws.Range("S17").Value = 0.1611649813509058
$ws.Range("T17").Value = 0.1611649813509058
$ws.Range("E18").Value = 2
$ws.Range("F18").Value = 0.6666666666666666
$ws.Range("G18").Value = 0.06234399999999999
$ws.Range("H18").Value = 0.187032
$ws.Range("I18").Value = 0.08319795590462073
$ws.Range("J18").Value = 0.08319795590462074
$ws.Range("M18").Value = 0.5396133333333334
$ws.Range("N18").Value = 1.61884
$ws.Range("O18").Value = 0.02587018426425635
$ws.Range("P18").Value = 0.02587018426425635
$ws.Range("Q18").Value = 0.03364165365333333
$ws.Range("R18").Value = 0.30277488288
$ws.Range("S18").Value = 0.002152346449662013
$ws.Range("T18").Value = 0.002152346449662014
$ws.Range("E19").Value = 2
$ws.Range("F19").Value = 0.6666666666666666
$ws.Range("G19").Value = 0.06234399999999999
$ws.Range("H19").Value = 0.187032
$ws.Range("I19").Value = 0.08319795590462073
$ws.Range("J19").Value = 0.08319795590462074
$ws.Range("M19").Value = 6.260434
$ws.Range("N19").Value = 18.781302
$ws.Range("O19").Value = 0.3001382122153186
$ws.Range("P19").Value = 0.3001382122153186
$ws.Range("Q19").Value = 0.3903004972959999
$ws.Range("R19").Value = 3.512704475664
$ws.Range("S19").Value = 0.02497088574518178
$ws.Range("T19").Value = 0.02497088574518178
$ws.Range("E20").Value = 2
$ws.Range("F20").Value = 0.6666666666666666
$ws.Range("G20").Value = 0.06234399999999999
$ws.Range("H20").Value = 0.187032
$ws.Range("I20").Value = 0.08319795590462073
$ws.Range("J20").Value = 0.08319795590462074
$ws.Range("M20").Value = 1.182122
$ws.Range("N20").Value = 3.546366
$ws.Range("O20").Value = 0.0566733845769154
$ws.Range("P20").Value = 0.0566733845769154
$ws.Range("Q20").Value = 0.07369821396799998
$ws.Range("R20").Value = 0.6632839257119999
$ws.Range("S20").Value = 0.004715109750995821
$ws.Range("T20").Value = 0.004715109750995821
$ws.Range("E21").Value = 2
$ws.Range("F21").Value = 0.6666666666666666
$ws.Range("G21").Value = 0.06234399999999999
$ws.Range("H21").Value = 0.187032
$ws.Range("I21").Value = 0.08319795590462073
$ws.Range("J21").Value = 0.08319795590462074
$ws.Range("K21").Value = 3
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 12.87633433333333
$ws.Range("N21").Value = 38.629003
$ws.Range("O21").Value = 0.6173182189435098
$ws.Range("P21").Value = 0.6173182189435098
$ws.Range("Q21").Value = 0.8027621876773331
$ws.Range("R21").Value = 7.224859689095998
$ws.Range("S21").Value = 0.05135961395878114
$ws.Range("T21").Value = 0.05135961395878114
